# AgileDeliveryGrowthPack.pptx update script
# - Bumps the cached "datetimeFigureOut" footer field from 03/11/2024 to 15/11/2024
#   on the slide master and every layout that carries a Date placeholder.
# - Renames "Catch Up" -> "Sync Up" (master title textbox) and
#   "Catch Up Template" -> "Sync Up Template" (slide 6 title).
# - Bumps "Author: Tom Hoyland | Version 0.3" -> "... Version 0.4" on every slide.

$p = $ppt.ActivePresentation

$oldDate = "03/11/2024"
$newDate = "15/11/2024"
$oldVersion = "Author: Tom Hoyland | Version 0.3"
$newVersion = "Author: Tom Hoyland | Version 0.4"

# ---- helper: walk a Shapes collection, update the first shape whose
#      PlaceholderFormat.Type is ppPlaceholderDate (16) and that currently
#      shows $oldDate; returns $true if it patched something.
function Update-DatePlaceholder($shapesCollection) {
    for ($di = 1; $di -le $shapesCollection.Count; $di++) {
        $shp = $shapesCollection.Item($di)
        $isDatePh = $false
        if ($shp.Type -eq 14) {
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDatePh = $true
                }
            } catch {
            }
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
                return $true
            }
        }
    }
    return $false
}

# ---- helper: walk a Shapes collection, update the first shape whose text
#      exactly equals $oldText to $newText; returns $true if patched.
function Update-ExactShapeText($shapesCollection, [string]$oldText, [string]$newText) {
    for ($ei = 1; $ei -le $shapesCollection.Count; $ei++) {
        $shp = $shapesCollection.Item($ei)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldText) {
                    $shp.TextFrame.TextRange.Text = $newText
                    return $true
                }
            }
        }
    }
    return $false
}

# ---- 1) Date placeholder on every Design's SlideMaster + all of its CustomLayouts.
for ($desi = 1; $desi -le $p.Designs.Count; $desi++) {
    $design = $p.Designs.Item($desi)
    $master = $design.SlideMaster

    [void](Update-DatePlaceholder $master.Shapes)

    for ($layi = 1; $layi -le $master.CustomLayouts.Count; $layi++) {
        $layout = $master.CustomLayouts.Item($layi)
        [void](Update-DatePlaceholder $layout.Shapes)
    }
}

# ---- 2) "Catch Up" -> "Sync Up" on the Custom Design's slide master textbox.
for ($desi2 = 1; $desi2 -le $p.Designs.Count; $desi2++) {
    $design2 = $p.Designs.Item($desi2)
    $master2 = $design2.SlideMaster
    [void](Update-ExactShapeText $master2.Shapes "Catch Up" "Sync Up")
}

# ---- 3) Per-slide updates: title "Catch Up Template" -> "Sync Up Template"
#         and footer "Author: ... Version 0.3" -> "... Version 0.4".
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    [void](Update-ExactShapeText $slide.Shapes "Catch Up Template" "Sync Up Template")
    [void](Update-ExactShapeText $slide.Shapes $oldVersion $newVersion)
}
